$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'25.857.04"
$ws.Range("E2").Value = "'  -0.17%  "

# Row 3
$ws.Range("D3").Value = "'1.629.85"
$ws.Range("E3").Value = "'  -0.66%  "

# Row 4
$ws.Range("E4").Value = "'  -0.01%  "

# Row 5
$ws.Range("D5").Value = "'213.89"
$ws.Range("E5").Value = "'  -0.54%  "

# Row 6
$ws.Range("D6").Value = "'0.5107"
$ws.Range("E6").Value = "'  +1.15%  "

# Row 7
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "'  -0.22%  "

# Row 8
$ws.Range("D8").Value = "'0.2545"
$ws.Range("E8").Value = "'  -1.12%  "

# Row 9
$ws.Range("D9").Value = "'0.06323"
$ws.Range("E9").Value = "'  -1.12%  "

# Row 10
$ws.Range("D10").Value = "'19.38"
$ws.Range("E10").Value = "'  -0.72%  "

# Row 11
$ws.Range("D11").Value = "'0.07741"
$ws.Range("E11").Value = "'  -0.49%  "

# Row 12
$ws.Range("D12").Value = "'4.264"
$ws.Range("E12").Value = "'  -0.01%  "

# Row 13
$ws.Range("D13").Value = "'1.639.89"
$ws.Range("E13").Value = "'  -0.34%  "

# Row 14
$ws.Range("D14").Value = "'0.5394"
$ws.Range("E14").Value = "'  -0.76%  "

# Row 15
$ws.Range("D15").Value = "'0.0₅7696"
$ws.Range("E15").Value = "'  -2.86%  "

# Row 16
$ws.Range("D16").Value = "'63.86"
$ws.Range("E16").Value = "'  -1.13%  "

# Row 17
$ws.Range("D17").Value = "'25.868.73"
$ws.Range("E17").Value = "'  -0.34%  "

# Row 18
$ws.Range("D18").Value = "'1.001"
$ws.Range("E18").Value = "'  -0.08%  "

# Row 19
$ws.Range("B19").Value = "'BitcoinCash"
$ws.Range("C19").Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "'194.39"
$ws.Range("E19").Value = "'  -1.95%  "

# Row 20
$ws.Range("B20").Value = "'Uniswap"
$ws.Range("C20").Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'4.410"
$ws.Range("E20").Value = "'  +0.88%  "

# Row 21
$ws.Range("D21").Value = "'9.875"
$ws.Range("E21").Value = "'  -0.09%  "

# Row 22
$ws.Range("D22").Value = "'5.999"
$ws.Range("E22").Value = "'  +0.49%  "

# Row 23
$ws.Range("D23").Value = "'1.003"
$ws.Range("E23").Value = "'  -0.18%  "

# Row 24
$ws.Range("D24").Value = "'1.855"
$ws.Range("E24").Value = "'  -0.82%  "

# Row 25
$ws.Range("D25").Value = "'140.51"
$ws.Range("E25").Value = "'  -0.51%  "

# Row 26
$ws.Range("D26").Value = "'0.1184"
$ws.Range("E26").Value = "'  +4.40%  "

# Row 27
$ws.Range("D27").Value = "'6.790"
$ws.Range("E27").Value = "'  -0.50%  "

# Row 28
$ws.Range("D28").Value = "'15.54"
$ws.Range("E28").Value = "'  -1.08%  "

# Row 29
$ws.Range("D29").Value = "'1.232"
$ws.Range("E29").Value = "'  -0.69%  "

# Row 30
$ws.Range("D30").Value = "'0.04893"
$ws.Range("E30").Value = "'  -0.90%  "

# Row 31
$ws.Range("D31").Value = "'3.234"
$ws.Range("E31").Value = "'  -0.92%  "

# Row 32
$ws.Range("D32").Value = "'3.150"
$ws.Range("E32").Value = "'  -1.35%  "

# Row 33
$ws.Range("D33").Value = "'1.520"
$ws.Range("E33").Value = "'  -0.89%  "

# Row 34
$ws.Range("D34").Value = "'2.362"
$ws.Range("E34").Value = "'  -0.09%  "

# Row 35
$ws.Range("D35").Value = "'0.8857"
$ws.Range("E35").Value = "'  -0.74%  "

# Row 36
$ws.Range("D36").Value = "'2.570"
$ws.Range("E36").Value = "'  -1.39%  "

# Row 37
$ws.Range("D37").Value = "'1.134.09"
$ws.Range("E37").Value = "'  -0.80%  "

# Row 38
$ws.Range("D38").Value = "'0.5377"
$ws.Range("E38").Value = "'  -3.23%  "

# Row 39
$ws.Range("D39").Value = "'0.01542"
$ws.Range("E39").Value = "'  -1.75%  "

# Row 40
$ws.Range("D40").Value = "'1.001"
$ws.Range("E40").Value = "'  -0.03%  "

# Row 41
$ws.Range("D41").Value = "'2.530"
$ws.Range("E41").Value = "'  -1.39%  "

# Row 42
$ws.Range("B42").Value = "'BabyDogeCoin"
$ws.Range("C42").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D42").Value = "'0.0₈126"
$ws.Range("E42").Value = "'  +5.24%  "

# Row 43
$ws.Range("B43").Value = "'TrustWalletToken"
$ws.Range("C43").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'0.8103"
$ws.Range("E43").Value = "'  +0.00%  "

# Row 44
$ws.Range("B44").Value = "'Quant"
$ws.Range("C44").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "'98.55"
$ws.Range("E44").Value = "'  -1.23%  "

# Row 45
$ws.Range("B45").Value = "'FraxShare"
$ws.Range("C45").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'5.429"
$ws.Range("E45").Value = "'  -4.83%  "

# Row 46
$ws.Range("D46").Value = "'1.768.40"
$ws.Range("E46").Value = "'  -0.58%  "

# Row 47
$ws.Range("D47").Value = "'0.4523"
$ws.Range("E47").Value = "'  +0.26%  "

# Row 48
$ws.Range("D48").Value = "'1.001"
$ws.Range("E48").Value = "'  -0.42%  "

# Row 49
$ws.Range("D49").Value = "'54.46"
$ws.Range("E49").Value = "'  -0.25%  "

# Row 50
$ws.Range("D50").Value = "'0.05047"
$ws.Range("E50").Value = "'  -0.18%  "

# Row 51
$ws.Range("D51").Value = "'1.002"
$ws.Range("E51").Value = "'  -0.27%  "
